$d = $word.ActiveDocument
$endPos = $d.Content.End
$rng = $d.Range($endPos, $endPos)
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>    <w:p/>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading2"/>
      </w:pPr>
      <w:r>
        <w:t>Repeating Stuff With Loops</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:tbl>
      <w:tblPr>
        <w:tblW w:w="0" w:type="auto"/>
        <w:tblCellSpacing w:w="15" w:type="dxa"/>
        <w:tblCellMar>
          <w:top w:w="15" w:type="dxa"/>
          <w:left w:w="15" w:type="dxa"/>
          <w:bottom w:w="15" w:type="dxa"/>
          <w:right w:w="15" w:type="dxa"/>
        </w:tblCellMar>
        <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="950"/>
        <w:gridCol w:w="2538"/>
        <w:gridCol w:w="2745"/>
      </w:tblGrid>
      <w:tr>
        <w:trPr>
          <w:tblHeader/>
          <w:tblCellSpacing w:w="15" w:type="dxa"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:vAlign w:val="center"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:lastRenderedPageBreak/>
              <w:t>Loop</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:vAlign w:val="center"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t>When to Use</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:vAlign w:val="center"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t>Best For</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:tblCellSpacing w:w="15" w:type="dxa"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:vAlign w:val="center"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>for</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:vAlign w:val="center"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Fixed number of repetitions</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:vAlign w:val="center"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Arrays, counters</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:tblCellSpacing w:w="15" w:type="dxa"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:vAlign w:val="center"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>while</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:vAlign w:val="center"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Unknown repetitions</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:vAlign w:val="center"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>User input, waiting conditions</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:tblCellSpacing w:w="15" w:type="dxa"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:vAlign w:val="center"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>do...while</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:vAlign w:val="center"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Must run once, then check</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:vAlign w:val="center"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Menu-like inputs</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:tblCellSpacing w:w="15" w:type="dxa"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:vAlign w:val="center"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>for...of</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:vAlign w:val="center"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Iterate values</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:vAlign w:val="center"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Arrays, strings</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:tblCellSpacing w:w="15" w:type="dxa"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:vAlign w:val="center"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>for...in</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:vAlign w:val="center"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Iterate keys</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:vAlign w:val="center"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Objects</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading2"/>
      </w:pPr>
      <w:r>
        <w:t>Looping over an arrays</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>c</w:t>
      </w:r>
      <w:r>
        <w:t>onst animals = [‘cheetah’, “Elephant”, “Dear”, “bear”]</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>f</w:t>
      </w:r>
      <w:r>
        <w:t>or(let i = 0; i &lt; animals.length; i++){</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:t>c</w:t>
      </w:r>
      <w:r>
        <w:t>onsole.log(animals[i]);</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading2"/>
      </w:pPr>
      <w:r>
        <w:t>Nested Loops</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>for(let i = 0; i &lt; 5; i++){</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">    console.log(`This is Outer Loop: ${i}`);</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">    for(let j = 0; j &lt; 2; j++){</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">        console.log("This is Inner LOOP");</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">    }</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>let myarray= [["ethishwar", "aadharaf", "aeeff"],</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">             ["Madacaraa", "chintu","Preethi"]];</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>console.log("Yethishwar");</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>for (let i = 0; i &lt; myarray.length; i++){</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">    console.log(myarray[i]);</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">    for(let j = 0; j &lt; myarray[i].length; j++){</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">        console.log(myarray[i][j]);</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">    }</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="CODE"/>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading2"/>
      </w:pPr>
      <w:r>
        <w:t>While Loop</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="9"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Generally used for fixed number of iterations</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading2"/>
      </w:pPr>
      <w:r>
        <w:t>Break keyword</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="9"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Used to stop the iteration after certain condition is met</w:t>
      </w:r>
    </w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
[void]$rng.InsertXML($xml)
Write-Output "InsertXML done"
Write-Output $d.Paragraphs.Count
Write-Output $d.Tables.Count
